$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column, using the same style as other headers (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New column data values (row 2 through row 10)
$saveValues = @(0, 1, 0, 0, 1, 1, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
